# Update countries & provincias Spain
# Daily refresh of the COVID table: new totals for several countries plus
# three countries (Angola/Haiti/Gabon) and one pair (Montserrat/Islas
# Malvinas) re-ranking in the (descending, by total cases) sorted list, and
# the "last updated" timestamp footer.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Datos actualizados ..." footer timestamp -----------------------------
$ws.Range("A1").Value = "Datos actualizados a 25 de Octubre de 2020 a las 02:35"

# --- Row 4: Estados Unidos --------------------------------------------------
$ws.Range("B4").Value = 8827154
$ws.Range("C4").Value = 78671
$ws.Range("D4").Value = 5741611
$ws.Range("E4").Value = 2855475
$ws.Range("G4").Value = 784
$ws.Range("H4").Value = 230068

# --- Row 112: Guayana Francesa ----------------------------------------------
$ws.Range("B112").Value = 10376
$ws.Range("C112").Value = 25
$ws.Range("E112").Value = 312

# --- Rows 113-117: Angola overtakes Haiti and Gabon in the ranking ---------
# Row 113 (Lituania) keeps its place.
# Row 114 now holds Angola (fresh data), row 115 now holds Haiti (its
# previous totals), row 116 now holds Gabon (its previous totals); row 117
# (Jamaica) is unaffected.
$ws.Range("A114").Value = "Angola"
$ws.Range("B114").Value = 9026
$ws.Range("C114").Value = 197
$ws.Range("D114").Value = 3461
$ws.Range("E114").Value = 5298
$ws.Range("G114").Value = 2
$ws.Range("H114").Value = 267

$ws.Range("A115").Value = "Haiti"
$ws.Range("B115").Value = 9015
$ws.Range("D115").Value = 7361
$ws.Range("E115").Value = 1423
$ws.Range("H115").Value = 231

$ws.Range("A116").Value = "Gabon"
$ws.Range("B116").Value = 8919
$ws.Range("D116").Value = 8512
$ws.Range("E116").Value = 353
$ws.Range("H116").Value = 54

# --- Row 135: Congo ----------------------------------------------------------
$ws.Range("B135").Value = 5253
$ws.Range("C135").Value = 97
$ws.Range("E135").Value = 1274

# --- Row 136: Surinam --------------------------------------------------------
$ws.Range("B136").Value = 5166
$ws.Range("C136").Value = 11
$ws.Range("D136").Value = 5016
$ws.Range("E136").Value = 41

# --- Row 142: Aruba -----------------------------------------------------------
$ws.Range("B142").Value = 4410
$ws.Range("C142").Value = 9
$ws.Range("D142").Value = 4186
$ws.Range("E142").Value = 188

# --- Row 147: Guyana -----------------------------------------------------------
$ws.Range("B147").Value = 3994
$ws.Range("C147").Value = 34
$ws.Range("D147").Value = 2970
$ws.Range("E147").Value = 907

# --- Row 161: Togo --------------------------------------------------------------
$ws.Range("B161").Value = 2187
$ws.Range("C161").Value = 25
$ws.Range("D161").Value = 1591
$ws.Range("E161").Value = 544

# --- Row 169: Santo Tome y Principe ---------------------------------------------
$ws.Range("B169").Value = 940
$ws.Range("C169").Value = 2
$ws.Range("D169").Value = 900

# --- Row 172: San Martin (Parte Holandesa) --------------------------------------
$ws.Range("B172").Value = 780
$ws.Range("C172").Value = 4
$ws.Range("D172").Value = 705
$ws.Range("E172").Value = 53

# --- Rows 216-217: Montserrat overtakes Islas Malvinas --------------------------
# Row 215 (San Pedro y Miquelon) and row 218 (Sahara Occidental) keep their place.
$ws.Range("A216").Value = "Montserrat"
$ws.Range("D216").Value = 12
$ws.Range("H216").Value = 1

$ws.Range("A217").Value = "Islas Malvinas"
$ws.Range("D217").Value = 13
$ws.Range("H217").Value = 0
